# Apply updated values to include "tri proximity" tables
$wb = $excel.ActiveWorkbook

# --- Sheet "Means" ---
$wsMeans = $wb.Worksheets.Item("Means")
$wsMeans.Range("D2").Value = 58
$wsMeans.Range("G2").Value = 51
$wsMeans.Range("D3").Value = 9.5
$wsMeans.Range("D5").Value = 67
$wsMeans.Range("D6").Value = 49
$wsMeans.Range("D7").Value = 12
$wsMeans.Range("G7").Value = 9.9
$wsMeans.Range("D8").Value = 9.9

# --- Sheet "Standard Deviations" ---
$wsSD = $wb.Worksheets.Item("Standard Deviations")
$wsSD.Range("D4").Value = 14
$wsSD.Range("D5").Value = 12
$wsSD.Range("D6").Value = 15
$wsSD.Range("D7").Value = 12
$wsSD.Range("D8").Value = 7.4
$wsSD.Range("G8").Value = 6.8
$wsSD.Range("D10").Value = 0.046
$wsSD.Range("F10").Value = 0.04
